# Excel COM-interop script implementing:
# "Finished Legionella positivity trend graph by facility area by source,
#  cleaned up all trend graphs with better colors, fontstyle and font sizing"
#
# Data-level effect on the "Analytical_Data" worksheet:
#   - Corrects several Legionella (col L) / Lp_Positive (col M) values for the
#     existing two sampling dates (rows 6,8,9,14,15,16,19,21,23,25).
#   - Appends two more full sampling-date blocks (rows 26-37 for 9/22/2021,
#     rows 38-49 for 12/24/2021) so the trend graphs have more data points.
#   - Updates the view (scroll position / active cell) to where the user was
#     last working while finishing the graphs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Legionella (L) / Lp_Positive (M) values on existing rows ---
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = "YES"
$ws.Cells.Item(8, 12).Value = 6
$ws.Cells.Item(8, 13).Value = "YES"
$ws.Cells.Item(9, 12).Value = 2
$ws.Cells.Item(9, 13).Value = "YES"
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = "YES"
$ws.Cells.Item(15, 12).Value = 60
$ws.Cells.Item(15, 13).Value = "YES"
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = "NO"
$ws.Cells.Item(19, 12).Value = 4
$ws.Cells.Item(19, 13).Value = "YES"
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = "YES"
$ws.Cells.Item(23, 12).Value = 3
$ws.Cells.Item(23, 13).Value = "YES"
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = "NO"

# --- Append new data rows 26-49 (two additional sampling dates) ---
# Row 26
$ws.Cells.Item(26, 1).Value = 44461
$ws.Cells.Item(26, 2).Value = "14:00:00"
$ws.Cells.Item(26, 3).Value = "CW"
$ws.Cells.Item(26, 4).Value = "MOB"
$ws.Cells.Item(26, 5).Value = 1
$ws.Cells.Item(26, 6).Value = "Supply"
$ws.Cells.Item(26, 7).Value = "Cooling Tower"
$ws.Cells.Item(26, 8).Value = "KMC"
$ws.Cells.Item(26, 9).Value = 78
$ws.Cells.Item(26, 10).Value = 0.04
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).Value = "NO"
$ws.Cells.Item(26, 14).Value = 1234

# Row 27
$ws.Cells.Item(27, 1).Value = 44461
$ws.Cells.Item(27, 2).Value = "14:00:00"
$ws.Cells.Item(27, 3).Value = "CW"
$ws.Cells.Item(27, 4).Value = "CP"
$ws.Cells.Item(27, 5).Value = 1
$ws.Cells.Item(27, 6).Value = "CT"
$ws.Cells.Item(27, 7).Value = "Cooling Tower"
$ws.Cells.Item(27, 8).Value = "KMC"
$ws.Cells.Item(27, 9).Value = 79
$ws.Cells.Item(27, 10).Value = 1.2
$ws.Cells.Item(27, 12).Value = 60
$ws.Cells.Item(27, 13).Value = "YES"
$ws.Cells.Item(27, 14).Value = 1235

# Row 28
$ws.Cells.Item(28, 1).Value = 44461
$ws.Cells.Item(28, 2).Value = "14:00:00"
$ws.Cells.Item(28, 3).Value = "PWC"
$ws.Cells.Item(28, 4).Value = "OLD"
$ws.Cells.Item(28, 5).Value = 2
$ws.Cells.Item(28, 6).Value = "SNK"
$ws.Cells.Item(28, 7).Value = "W241"
$ws.Cells.Item(28, 8).Value = "KMC"
$ws.Cells.Item(28, 9).Value = 72
$ws.Cells.Item(28, 11).Value = 0.8
$ws.Cells.Item(28, 12).Value = 2
$ws.Cells.Item(28, 13).Value = "YES"
$ws.Cells.Item(28, 14).Value = 1236

# Row 29
$ws.Cells.Item(29, 1).Value = 44461
$ws.Cells.Item(29, 2).Value = "14:00:00"
$ws.Cells.Item(29, 3).Value = "PWC"
$ws.Cells.Item(29, 4).Value = "NEW"
$ws.Cells.Item(29, 5).Value = 1
$ws.Cells.Item(29, 6).Value = "SNK"
$ws.Cells.Item(29, 7).Value = "CCU ROOM 1"
$ws.Cells.Item(29, 8).Value = "KMC"
$ws.Cells.Item(29, 9).Value = 69
$ws.Cells.Item(29, 11).Value = 1.4
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = "NO"
$ws.Cells.Item(29, 14).Value = 1237

# Row 30
$ws.Cells.Item(30, 1).Value = 44461
$ws.Cells.Item(30, 2).Value = "14:00:00"
$ws.Cells.Item(30, 3).Value = "PWC"
$ws.Cells.Item(30, 4).Value = "NEW"
$ws.Cells.Item(30, 5).Value = 2
$ws.Cells.Item(30, 6).Value = "SNK"
$ws.Cells.Item(30, 7).Value = 254
$ws.Cells.Item(30, 8).Value = "KMC"
$ws.Cells.Item(30, 9).Value = 70
$ws.Cells.Item(30, 11).Value = 1.2
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = "NO"
$ws.Cells.Item(30, 14).Value = 1238

# Row 31
$ws.Cells.Item(31, 1).Value = 44461
$ws.Cells.Item(31, 2).Value = "14:00:00"
$ws.Cells.Item(31, 3).Value = "PWH"
$ws.Cells.Item(31, 4).Value = "NEW"
$ws.Cells.Item(31, 5).Value = 2
$ws.Cells.Item(31, 6).Value = "SNK"
$ws.Cells.Item(31, 7).Value = 254
$ws.Cells.Item(31, 8).Value = "KMC"
$ws.Cells.Item(31, 9).Value = 114
$ws.Cells.Item(31, 11).Value = 0.2
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = "NO"
$ws.Cells.Item(31, 14).Value = 1239

# Row 32
$ws.Cells.Item(32, 1).Value = 44461
$ws.Cells.Item(32, 2).Value = "14:00:00"
$ws.Cells.Item(32, 3).Value = "PWH"
$ws.Cells.Item(32, 4).Value = "NEW"
$ws.Cells.Item(32, 5).Value = 2
$ws.Cells.Item(32, 6).Value = "SNK"
$ws.Cells.Item(32, 7).Value = 259
$ws.Cells.Item(32, 8).Value = "KMC"
$ws.Cells.Item(32, 9).Value = 116
$ws.Cells.Item(32, 11).Value = 0.4
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = "NO"
$ws.Cells.Item(32, 14).Value = 1240

# Row 33
$ws.Cells.Item(33, 1).Value = 44461
$ws.Cells.Item(33, 2).Value = "14:00:00"
$ws.Cells.Item(33, 3).Value = "PWH"
$ws.Cells.Item(33, 4).Value = "NEW"
$ws.Cells.Item(33, 5).Value = 1
$ws.Cells.Item(33, 6).Value = "SNK"
$ws.Cells.Item(33, 7).Value = "CCU ROOM 1"
$ws.Cells.Item(33, 8).Value = "KMC"
$ws.Cells.Item(33, 9).Value = 118
$ws.Cells.Item(33, 11).Value = 0.3
$ws.Cells.Item(33, 12).Value = 1
$ws.Cells.Item(33, 13).Value = "YES"
$ws.Cells.Item(33, 14).Value = 1241

# Row 34
$ws.Cells.Item(34, 1).Value = 44461
$ws.Cells.Item(34, 2).Value = "14:00:00"
$ws.Cells.Item(34, 3).Value = "PWH"
$ws.Cells.Item(34, 4).Value = "OLD"
$ws.Cells.Item(34, 5).Value = 2
$ws.Cells.Item(34, 6).Value = "SNK"
$ws.Cells.Item(34, 7).Value = "2W241"
$ws.Cells.Item(34, 8).Value = "KMC"
$ws.Cells.Item(34, 9).Value = 108
$ws.Cells.Item(34, 11).Value = 0.08
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = "NO"
$ws.Cells.Item(34, 14).Value = 1242

# Row 35
$ws.Cells.Item(35, 1).Value = 44461
$ws.Cells.Item(35, 2).Value = "14:00:00"
$ws.Cells.Item(35, 3).Value = "PWH"
$ws.Cells.Item(35, 4).Value = "OLD"
$ws.Cells.Item(35, 5).Value = 2
$ws.Cells.Item(35, 6).Value = "SNK"
$ws.Cells.Item(35, 7).Value = "2W NURSE PANTRY ACROSS FROM 234"
$ws.Cells.Item(35, 8).Value = "KMC"
$ws.Cells.Item(35, 9).Value = 117
$ws.Cells.Item(35, 11).Value = 0.6
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = "NO"
$ws.Cells.Item(35, 14).Value = 1243

# Row 36
$ws.Cells.Item(36, 1).Value = 44461
$ws.Cells.Item(36, 2).Value = "14:00:00"
$ws.Cells.Item(36, 3).Value = "PWH"
$ws.Cells.Item(36, 4).Value = "OLD"
$ws.Cells.Item(36, 5).Value = 2
$ws.Cells.Item(36, 6).Value = "SNK"
$ws.Cells.Item(36, 7).Value = "E202"
$ws.Cells.Item(36, 8).Value = "KMC"
$ws.Cells.Item(36, 9).Value = 116
$ws.Cells.Item(36, 11).Value = 0.3
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = "NO"
$ws.Cells.Item(36, 14).Value = 1244

# Row 37
$ws.Cells.Item(37, 1).Value = 44461
$ws.Cells.Item(37, 2).Value = "14:00:00"
$ws.Cells.Item(37, 3).Value = "IC"
$ws.Cells.Item(37, 4).Value = "OLD"
$ws.Cells.Item(37, 5).Value = 2
$ws.Cells.Item(37, 6).Value = "IM"
$ws.Cells.Item(37, 7).Value = "ICE IM07 2W NURSE PANTRY ACROSS FROM 234"
$ws.Cells.Item(37, 8).Value = "KMC"
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 13).Value = "NO"
$ws.Cells.Item(37, 14).Value = 1245

# Row 38
$ws.Cells.Item(38, 1).Value = 44554
$ws.Cells.Item(38, 2).Value = "14:00:00"
$ws.Cells.Item(38, 3).Value = "CW"
$ws.Cells.Item(38, 4).Value = "MOB"
$ws.Cells.Item(38, 5).Value = 1
$ws.Cells.Item(38, 6).Value = "Supply"
$ws.Cells.Item(38, 7).Value = "Cooling Tower"
$ws.Cells.Item(38, 8).Value = "KMC"
$ws.Cells.Item(38, 9).Value = 78
$ws.Cells.Item(38, 10).Value = 0.04
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 13).Value = "NO"
$ws.Cells.Item(38, 14).Value = 1234

# Row 39
$ws.Cells.Item(39, 1).Value = 44554
$ws.Cells.Item(39, 2).Value = "14:00:00"
$ws.Cells.Item(39, 3).Value = "CW"
$ws.Cells.Item(39, 4).Value = "CP"
$ws.Cells.Item(39, 5).Value = 1
$ws.Cells.Item(39, 6).Value = "CT"
$ws.Cells.Item(39, 7).Value = "Cooling Tower"
$ws.Cells.Item(39, 8).Value = "KMC"
$ws.Cells.Item(39, 9).Value = 79
$ws.Cells.Item(39, 10).Value = 1.2
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).Value = "NO"
$ws.Cells.Item(39, 14).Value = 1235

# Row 40
$ws.Cells.Item(40, 1).Value = 44554
$ws.Cells.Item(40, 2).Value = "14:00:00"
$ws.Cells.Item(40, 3).Value = "PWC"
$ws.Cells.Item(40, 4).Value = "OLD"
$ws.Cells.Item(40, 5).Value = 2
$ws.Cells.Item(40, 6).Value = "SNK"
$ws.Cells.Item(40, 7).Value = "W241"
$ws.Cells.Item(40, 8).Value = "KMC"
$ws.Cells.Item(40, 9).Value = 72
$ws.Cells.Item(40, 11).Value = 0.8
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = "NO"
$ws.Cells.Item(40, 14).Value = 1236

# Row 41
$ws.Cells.Item(41, 1).Value = 44554
$ws.Cells.Item(41, 2).Value = "14:00:00"
$ws.Cells.Item(41, 3).Value = "PWC"
$ws.Cells.Item(41, 4).Value = "NEW"
$ws.Cells.Item(41, 5).Value = 1
$ws.Cells.Item(41, 6).Value = "SNK"
$ws.Cells.Item(41, 7).Value = "CCU ROOM 1"
$ws.Cells.Item(41, 8).Value = "KMC"
$ws.Cells.Item(41, 9).Value = 69
$ws.Cells.Item(41, 11).Value = 1.4
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).Value = "NO"
$ws.Cells.Item(41, 14).Value = 1237

# Row 42
$ws.Cells.Item(42, 1).Value = 44554
$ws.Cells.Item(42, 2).Value = "14:00:00"
$ws.Cells.Item(42, 3).Value = "PWC"
$ws.Cells.Item(42, 4).Value = "NEW"
$ws.Cells.Item(42, 5).Value = 2
$ws.Cells.Item(42, 6).Value = "SNK"
$ws.Cells.Item(42, 7).Value = 254
$ws.Cells.Item(42, 8).Value = "KMC"
$ws.Cells.Item(42, 9).Value = 70
$ws.Cells.Item(42, 11).Value = 1.2
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = "NO"
$ws.Cells.Item(42, 14).Value = 1238

# Row 43
$ws.Cells.Item(43, 1).Value = 44554
$ws.Cells.Item(43, 2).Value = "14:00:00"
$ws.Cells.Item(43, 3).Value = "PWH"
$ws.Cells.Item(43, 4).Value = "NEW"
$ws.Cells.Item(43, 5).Value = 2
$ws.Cells.Item(43, 6).Value = "SNK"
$ws.Cells.Item(43, 7).Value = 254
$ws.Cells.Item(43, 8).Value = "KMC"
$ws.Cells.Item(43, 9).Value = 114
$ws.Cells.Item(43, 11).Value = 0.2
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = "NO"
$ws.Cells.Item(43, 14).Value = 1239

# Row 44
$ws.Cells.Item(44, 1).Value = 44554
$ws.Cells.Item(44, 2).Value = "14:00:00"
$ws.Cells.Item(44, 3).Value = "PWH"
$ws.Cells.Item(44, 4).Value = "NEW"
$ws.Cells.Item(44, 5).Value = 2
$ws.Cells.Item(44, 6).Value = "SNK"
$ws.Cells.Item(44, 7).Value = 259
$ws.Cells.Item(44, 8).Value = "KMC"
$ws.Cells.Item(44, 9).Value = 116
$ws.Cells.Item(44, 11).Value = 0.4
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).Value = "NO"
$ws.Cells.Item(44, 14).Value = 1240

# Row 45
$ws.Cells.Item(45, 1).Value = 44554
$ws.Cells.Item(45, 2).Value = "14:00:00"
$ws.Cells.Item(45, 3).Value = "PWH"
$ws.Cells.Item(45, 4).Value = "NEW"
$ws.Cells.Item(45, 5).Value = 1
$ws.Cells.Item(45, 6).Value = "SNK"
$ws.Cells.Item(45, 7).Value = "CCU ROOM 1"
$ws.Cells.Item(45, 8).Value = "KMC"
$ws.Cells.Item(45, 9).Value = 118
$ws.Cells.Item(45, 11).Value = 0.3
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = "NO"
$ws.Cells.Item(45, 14).Value = 1241

# Row 46
$ws.Cells.Item(46, 1).Value = 44554
$ws.Cells.Item(46, 2).Value = "14:00:00"
$ws.Cells.Item(46, 3).Value = "PWH"
$ws.Cells.Item(46, 4).Value = "OLD"
$ws.Cells.Item(46, 5).Value = 2
$ws.Cells.Item(46, 6).Value = "SNK"
$ws.Cells.Item(46, 7).Value = "2W241"
$ws.Cells.Item(46, 8).Value = "KMC"
$ws.Cells.Item(46, 9).Value = 108
$ws.Cells.Item(46, 11).Value = 0.08
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = "NO"
$ws.Cells.Item(46, 14).Value = 1242

# Row 47
$ws.Cells.Item(47, 1).Value = 44554
$ws.Cells.Item(47, 2).Value = "14:00:00"
$ws.Cells.Item(47, 3).Value = "PWH"
$ws.Cells.Item(47, 4).Value = "OLD"
$ws.Cells.Item(47, 5).Value = 2
$ws.Cells.Item(47, 6).Value = "SNK"
$ws.Cells.Item(47, 7).Value = "2W NURSE PANTRY ACROSS FROM 234"
$ws.Cells.Item(47, 8).Value = "KMC"
$ws.Cells.Item(47, 9).Value = 117
$ws.Cells.Item(47, 11).Value = 0.6
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = "NO"
$ws.Cells.Item(47, 14).Value = 1243

# Row 48
$ws.Cells.Item(48, 1).Value = 44554
$ws.Cells.Item(48, 2).Value = "14:00:00"
$ws.Cells.Item(48, 3).Value = "PWH"
$ws.Cells.Item(48, 4).Value = "OLD"
$ws.Cells.Item(48, 5).Value = 2
$ws.Cells.Item(48, 6).Value = "SNK"
$ws.Cells.Item(48, 7).Value = "E202"
$ws.Cells.Item(48, 8).Value = "KMC"
$ws.Cells.Item(48, 9).Value = 116
$ws.Cells.Item(48, 11).Value = 0.3
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 13).Value = "NO"
$ws.Cells.Item(48, 14).Value = 1244

# Row 49
$ws.Cells.Item(49, 1).Value = 44554
$ws.Cells.Item(49, 2).Value = "14:00:00"
$ws.Cells.Item(49, 3).Value = "IC"
$ws.Cells.Item(49, 4).Value = "OLD"
$ws.Cells.Item(49, 5).Value = 2
$ws.Cells.Item(49, 6).Value = "IM"
$ws.Cells.Item(49, 7).Value = "ICE IM07 2W NURSE PANTRY ACROSS FROM 234"
$ws.Cells.Item(49, 8).Value = "KMC"
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 13).Value = "NO"
$ws.Cells.Item(49, 14).Value = 1245

# --- Update the sheet view: scroll down toward the newly added rows and
#     leave the active selection where the author left off (cell O33). ---
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 11
} catch {
}
$ws.Range("O33").Select()
